$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2097.93
$ws.Range("J17").Value = 2097.93
$ws.Range("L17").Value = 6293.789999999999
$ws.Range("N17").Value = -6629.789999999999
$ws.Range("H28").Value = 68766.8
$ws.Range("I28").Value = 79266.46000000001
$ws.Range("K28").Value = 79266.46000000001
$ws.Range("M28").Value = -78781.46000000001
$ws.Range("H112").Value = 1689.3
$ws.Range("J112").Value = 1720.3158
$ws.Range("L112").Value = 5160.9474
$ws.Range("N112").Value = -7376.9474
$ws.Range("H125").Value = 10104193
$ws.Range("I125").Value = 2108.5
$ws.Range("J125").Value = 12349101
$ws.Range("K125").Value = 18976.5
$ws.Range("L125").Value = 111141909
$ws.Range("M125").Value = -16516.5
$ws.Range("N125").Value = -111146829
$ws.Range("H132").Value = 2181.6667
$ws.Range("I132").Value = 2132.3948
$ws.Range("K132").Value = 6397.1844
$ws.Range("M132").Value = -3867.1844
$ws.Range("H135").Value = 2091.7188
$ws.Range("I135").Value = 1877.3103
$ws.Range("J135").Value = 4164.3335
$ws.Range("K135").Value = 16895.7927
$ws.Range("L135").Value = 37479.0015
$ws.Range("M135").Value = -14360.7927
$ws.Range("N135").Value = -42549.0015
$ws.Range("H138").Value = 5883
$ws.Range("I138").Value = 3627.3809
$ws.Range("J138").Value = 6776.736
$ws.Range("K138").Value = 10882.1427
$ws.Range("L138").Value = 20330.208
$ws.Range("M138").Value = -5742.1427
$ws.Range("N138").Value = -30610.208

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6751.5557
$ws.Range("I2").Value = 7979.4287
$ws.Range("K2").Value = 7979.4287
$ws.Range("M2").Value = -7866.4287
$ws.Range("H32").Value = 4743.732
$ws.Range("I32").Value = 4125.4907
$ws.Range("K32").Value = 4125.4907
$ws.Range("M32").Value = -3838.4907
$ws.Range("H61").Value = 3442.889
$ws.Range("I61").Value = 3251.0688
$ws.Range("J61").Value = 4237.5713
$ws.Range("K61").Value = 3251.0688
$ws.Range("L61").Value = 4237.5713
$ws.Range("M61").Value = -3039.0688
$ws.Range("N61").Value = -4661.5713
$ws.Range("H110").Value = 180573
$ws.Range("I110").Value = 229385.73
$ws.Range("K110").Value = 229385.73
$ws.Range("M110").Value = -227340.73
$ws.Range("H116").Value = 6751.5557
$ws.Range("I116").Value = 7979.4287
$ws.Range("K116").Value = 7979.4287
$ws.Range("M116").Value = -5685.4287
$ws.Range("H132").Value = 3150.054
$ws.Range("J132").Value = 5259.524
$ws.Range("L132").Value = 15778.572
$ws.Range("N132").Value = -20838.572
$ws.Range("H136").Value = 3442.889
$ws.Range("I136").Value = 3251.0688
$ws.Range("J136").Value = 4237.5713
$ws.Range("K136").Value = 9753.206399999999
$ws.Range("L136").Value = 12712.7139
$ws.Range("M136").Value = -7203.206399999999
$ws.Range("N136").Value = -17812.7139
$ws.Range("H139").Value = 49571.668
$ws.Range("J139").Value = 49571.668
$ws.Range("L139").Value = 49571.668
$ws.Range("N139").Value = -59851.668

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6751.5557
$ws.Range("I3").Value = 7979.4287
$ws.Range("K3").Value = 7979.4287
$ws.Range("M3").Value = -7865.4287
$ws.Range("H80").Value = 156.63637
$ws.Range("J80").Value = 177.42857
$ws.Range("L80").Value = 177.42857
$ws.Range("N80").Value = -2173.42857
$ws.Range("H83").Value = 156.63637
$ws.Range("J83").Value = 177.42857
$ws.Range("L83").Value = 887.1428500000001
$ws.Range("N83").Value = -10871.14285
$ws.Range("H107").Value = 371854.22
$ws.Range("I107").Value = 1268.3158
$ws.Range("J107").Value = 1251995.8
$ws.Range("K107").Value = 1268.3158
$ws.Range("L107").Value = 1251995.8
$ws.Range("M107").Value = 651.6841999999999
$ws.Range("N107").Value = -1255835.8

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4363.037
$ws.Range("I16").Value = 3904.9443
$ws.Range("K16").Value = 3904.9443
$ws.Range("M16").Value = -3617.9443
$ws.Range("H55").Value = 10000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 10000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 10000
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -10630
$ws.Range("H94").Value = 1739.8
$ws.Range("J94").Value = 2266.3333
$ws.Range("L94").Value = 2266.3333
$ws.Range("N94").Value = -3168.3333
$ws.Range("H105").Value = 1112.3334
$ws.Range("I105").Value = 1068.5
$ws.Range("K105").Value = 1068.5
$ws.Range("M105").Value = 678.5
$ws.Range("H113").Value = 4363.037
$ws.Range("I113").Value = 3904.9443
$ws.Range("K113").Value = 3904.9443
$ws.Range("M113").Value = -1734.9443
$ws.Range("H132").Value = 1471.62
$ws.Range("I132").Value = 1256.4889
$ws.Range("K132").Value = 3769.4667
$ws.Range("M132").Value = -1239.4667

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 5625330
$ws.Range("I32").Value = 30000000
$ws.Range("J32").Value = 2143234.5
$ws.Range("K32").Value = 90000000
$ws.Range("L32").Value = 6429703.5
$ws.Range("M32").Value = -89999717
$ws.Range("N32").Value = -6430269.5
$ws.Range("H56").Value = 5708.5454
$ws.Range("I56").Value = 5708.5454
$ws.Range("K56").Value = 5708.5454
$ws.Range("M56").Value = -5178.5454
$ws.Range("H137").Value = 3284.9443
$ws.Range("I137").Value = 3019.8
$ws.Range("J137").Value = 4610.6665
$ws.Range("K137").Value = 9059.400000000001
$ws.Range("L137").Value = 13831.9995
$ws.Range("M137").Value = -3959.400000000001
$ws.Range("N137").Value = -24031.9995

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 921769.0600000001
$ws.Range("I113").Value = 3335336.8
$ws.Range("K113").Value = 3335336.8
$ws.Range("M113").Value = -3333166.8

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H132").Value = 2547.8333
$ws.Range("I132").Value = 2236.1482
$ws.Range("K132").Value = 6708.444600000001
$ws.Range("M132").Value = -4178.444600000001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 18000
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H100").Value = 1445.9286
$ws.Range("I100").Value = 1678.9
$ws.Range("K100").Value = 3357.8
$ws.Range("M100").Value = -2816.8
$ws.Range("H132").Value = 66217
$ws.Range("I132").Value = 9223
$ws.Range("J132").Value = 116878.336
$ws.Range("K132").Value = 27669
$ws.Range("L132").Value = 350635.008
$ws.Range("M132").Value = -25139
$ws.Range("N132").Value = -355695.008
$ws.Range("H136").Value = 43731.293
$ws.Range("I136").Value = 10154.451
$ws.Range("J136").Value = 288362.56
$ws.Range("K136").Value = 30463.353
$ws.Range("L136").Value = 865087.6799999999
$ws.Range("M136").Value = -27913.353
$ws.Range("N136").Value = -870187.6799999999
